$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Include from Reasons why cert")

# --- Metadata sheet edits ---

# URL: pythia -> cicada
$ws1.Range("B2").Value = "http://fhirfli.dev/fhir/ig/cicada/ValueSet/preferred-allowed-reason"

# Date: updated timestamp
$ws1.Range("B8").Value = "2026-02-11T14:37:07-05:00"

# Insert a new "Jurisdiction" row after "Contact" (row 10), shifting rows 11-14 down to 12-15.
# Use Copy(destination) (bottom-up) so the existing cell style/formatting (s="2") carries over,
# then fix up values explicitly (Copy does not blank a populated destination when the source is empty).
$ws1.Range("A14:B14").Copy($ws1.Range("A15:B15"))
$ws1.Range("A13:B13").Copy($ws1.Range("A14:B14"))
$ws1.Range("A12:B12").Copy($ws1.Range("A13:B13"))
$ws1.Range("A11:B11").Copy($ws1.Range("A12:B12"))
$excel.CutCopyMode = $false

$ws1.Range("B14").ClearContents()
$ws1.Range("B13").ClearContents()

$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").ClearContents()

# --- Include sheet edits ---

# System URI: pythia -> cicada
$ws2.Range("B9").Value = "http://fhirfli.dev/fhir/ig/cicada/CodeSystem/PreferredAllowedReason"

# Rename the Include sheet
$ws2.Name = "Include #0"
